# Apply cryptos.xlsx price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal string that Excel would otherwise
# auto-convert to a number (losing e.g. trailing zeros) - force text format first.
$textCells = @("D4", "D5", "D6", "D11", "D17", "D18", "D19", "D20", "D21", "D24", "D29", "D30", "D36", "D37", "D38", "D39", "D40", "D41", "D44", "D45", "D46", "D48", "D49", "D51")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '64.139.86'
$ws.Range("E2").Value = '  -0.20%  '

$ws.Range("D3").Value = '3.478.72'
$ws.Range("E3").Value = '  -0.65%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = '585.75'
$ws.Range("E5").Value = '  -0.12%  '

$ws.Range("D6").Value = '132.08'
$ws.Range("E6").Value = '  -1.66%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").Value = '  -0.54%  '

$ws.Range("E9").Value = '  +5.04%  '

$ws.Range("E10").Value = '  -1.65%  '

$ws.Range("D11").Value = '0.386'
$ws.Range("E11").Value = '  +0.00%  '

$ws.Range("D12").Value = '4.074.79'
$ws.Range("E12").Value = '  -0.57%  '

$ws.Range("E13").Value = '  +0.02%  '

$ws.Range("E14").Value = '  -2.10%  '

$ws.Range("D15").Value = '3.482.16'
$ws.Range("E15").Value = '  -0.55%  '

$ws.Range("D16").Value = '64.112.23'
$ws.Range("E16").Value = '  -0.31%  '

$ws.Range("D17").Value = '24.28'
$ws.Range("E17").Value = '  -7.16%  '

$ws.Range("D18").Value = '9.96'
$ws.Range("E18").Value = '  +0.43%  '

$ws.Range("D19").Value = '5.72'
$ws.Range("E19").Value = '  -0.45%  '

$ws.Range("D20").Value = '13.49'
$ws.Range("E20").Value = '  -2.12%  '

$ws.Range("D21").Value = '384.87'
$ws.Range("E21").Value = '  -2.13%  '

$ws.Range("E22").Value = '  +0.55%  '

$ws.Range("D23").Value = '3.619.47'
$ws.Range("E23").Value = '  -0.61%  '

$ws.Range("D24").Value = '74.73'
$ws.Range("E24").Value = '  +0.71%  '

$ws.Range("E25").Value = '  -0.09%  '

$ws.Range("E26").Value = '  +0.08%  '

$ws.Range("E27").Value = '  -1.75%  '

$ws.Range("E28").Value = '  +0.15%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '7.18'
$ws.Range("E30").Value = '  -4.62%  '

$ws.Range("E31").Value = '  -4.87%  '

$ws.Range("E32").Value = '  -4.06%  '

$ws.Range("E33").Value = '  +2.72%  '

$ws.Range("D34").Value = '3.509.68'
$ws.Range("E34").Value = '  -0.39%  '

$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("D36").Value = '22.97'
$ws.Range("E36").Value = '  -2.19%  '

$ws.Range("D37").Value = '5.19'
$ws.Range("E37").Value = '  -0.43%  '

$ws.Range("D38").Value = '6.78'
$ws.Range("E38").Value = '  -2.01%  '

$ws.Range("D39").Value = '1.51'
$ws.Range("E39").Value = '  -3.58%  '

$ws.Range("D40").Value = '163.10'
$ws.Range("E40").Value = '  +0.69%  '

$ws.Range("D41").Value = '0.0776'
$ws.Range("E41").Value = '  -0.76%  '

$ws.Range("E42").Value = '  -0.72%  '

$ws.Range("E43").Value = '  +0.08%  '

$ws.Range("D44").Value = '4.32'
$ws.Range("E44").Value = '  -2.00%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '23.80'
$ws.Range("E45").Value = '  -6.50%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = '1.63'
$ws.Range("E46").Value = '  -1.62%  '

$ws.Range("E47").Value = '  -3.08%  '

$ws.Range("D48").Value = '0.916'
$ws.Range("E48").Value = '  +2.32%  '

$ws.Range("D49").Value = '6.71'
$ws.Range("E49").Value = '  -1.27%  '

$ws.Range("D50").Value = '2.358.79'
$ws.Range("E50").Value = '  -4.53%  '

$ws.Range("D51").Value = '0.0254'
$ws.Range("E51").Value = '  -3.00%  '
